# Update the "want to go" (想去人数) counts in column F on the
# "展览" (Exhibition) and "全部类型" (All types) worksheets.
# These two sheets contain duplicate rows for the same events, so the
# same new values are applied to both.

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> list of (row, newValue) pairs to update in column F.
$updates = @{
    "展览" = @(
        @{ Row = 2;  Value = 7211 },
        @{ Row = 7;  Value = 181 },
        @{ Row = 8;  Value = 130 },
        @{ Row = 12; Value = 217 },
        @{ Row = 13; Value = 13 },
        @{ Row = 19; Value = 3752 },
        @{ Row = 26; Value = 2423 },
        @{ Row = 36; Value = 27 },
        @{ Row = 37; Value = 165 },
        @{ Row = 38; Value = 1444 },
        @{ Row = 39; Value = 148 }
    )
    "全部类型" = @(
        @{ Row = 2;  Value = 7211 },
        @{ Row = 8;  Value = 181 },
        @{ Row = 9;  Value = 130 },
        @{ Row = 13; Value = 217 },
        @{ Row = 14; Value = 13 },
        @{ Row = 20; Value = 3752 },
        @{ Row = 27; Value = 2423 },
        @{ Row = 37; Value = 27 },
        @{ Row = 38; Value = 165 },
        @{ Row = 39; Value = 1444 },
        @{ Row = 40; Value = 148 }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $updates[$sheetName]) {
        $ws.Cells.Item($entry.Row, 6).Value = $entry.Value
    }
}
